$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheet name:" $ws.Name
Write-Host "B2 value:" $ws.Range("B2").Value
